# --- sheet2 ("Stage and Bosses"): world 2 / world 3 stage+boss table rewrite ---
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Basic enemies"
$ws2 = $wb.Worksheets.Item(2)   # "Stage and Bosses"

# 1) Re-point cell formats BEFORE overwriting values, copying from rows that
#    already carry the target fill colour so the colour-coded row banding
#    (one colour per "world") stays correct after rows are re-shuffled.
$ws2.Range("A9:E9").Copy()
$ws2.Range("A10:E10").PasteSpecial(-4122)   # xlPasteFormats

$ws2.Range("A17:E17").Copy()
$ws2.Range("A18:E18").PasteSpecial(-4122)

$ws2.Range("A16:E16").Copy()
$ws2.Range("A17:E17").PasteSpecial(-4122)

$ws2.Range("A19:E19").Copy()
$ws2.Range("A20:E20").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# 2) Clear the two stray "basic boss" gimmick tags on the Ultimatebasic/Carrier rows.
$ws2.Range("C2").ClearContents()
$ws2.Range("C3").ClearContents()

# 3) Rewrite the Stage / Name / Gimmick columns for rows 2-19 with the fleshed-out
#    world 2 & world 3 level + boss design notes, and add the new secret-boss row 20.
$ws2.Cells.Item(2,1).Value2 = "[1-10]"
$ws2.Cells.Item(2,2).Value2 = "Ultimatebasic"
$ws2.Cells.Item(3,1).Value2 = "[1-20]"
$ws2.Cells.Item(3,2).Value2 = "Carrier"
$ws2.Cells.Item(4,1).Value2 = "[1-25]"
$ws2.Cells.Item(4,2).Value2 = "Core"
$ws2.Cells.Item(5,1).Value2 = "[2-10]"
$ws2.Cells.Item(5,2).Value2 = "mix ultimates, buff, debuff"
$ws2.Cells.Item(6,1).Value2 = "[2-15]"
$ws2.Cells.Item(6,2).Value2 = "mesocarrier"
$ws2.Cells.Item(7,1).Value2 = "[2-18]"
$ws2.Cells.Item(7,2).Value2 = "Couplad"
$ws2.Cells.Item(7,3).Value2 = "2 entities are coupled to each other. When 1 dies, it slowly regains hp. If it reaches full hp then it is fully functional again, the only way to kill it is to the kill the other entitiy while one is down. (deals DOT)"
$ws2.Cells.Item(8,1).Value2 = "[2-20]"
$ws2.Cells.Item(8,2).Value2 = "mesocore"
$ws2.Cells.Item(9,1).Value2 = "[2-25]"
$ws2.Cells.Item(9,2).Value2 = "meso carrier + core"
$ws2.Cells.Item(10,1).Value2 = "[2-30]"
$ws2.Cells.Item(10,2).Value2 = "Minima"
$ws2.Cells.Item(10,3).Value2 = "small hitbox, buffs speed, moves towards randomm positions in the map, phases in and out of invisibility under 50% health & speed boost."
$ws2.Cells.Item(11,1).Value2 = "[3-10]"
$ws2.Cells.Item(11,2).Value2 = "MaxCouplad"
$ws2.Cells.Item(11,3).Value2 = "same as couplad but now everytime 1 of the entities dies it boost a certain aspect, 1 of the sides boost the speed and other boost damage (DOT increases with each kill as well)"
$ws2.Cells.Item(12,1).Value2 = "[3-15]"
$ws2.Cells.Item(12,2).Value2 = "hypercarrier"
$ws2.Cells.Item(13,1).Value2 = "[3-20]"
$ws2.Cells.Item(13,2).Value2 = "hypercore"
$ws2.Cells.Item(14,1).Value2 = "[3-25]"
$ws2.Cells.Item(14,2).Value2 = "Gigantodon"
$ws2.Cells.Item(14,3).Value2 = "masssive hp, very wide (pierce to hit stuff behind it), under 20%health recovers hp. Increases your ammo rate but makes it 0 every minute."
$ws2.Cells.Item(15,1).Value2 = "[3-35]"
$ws2.Cells.Item(15,2).Value2 = "coupled, maxcouplad"
$ws2.Cells.Item(16,1).Value2 = "[3-40]"
$ws2.Cells.Item(16,2).Value2 = "Minima & Gigantodon"
$ws2.Cells.Item(17,1).Value2 = "[3-45]"
$ws2.Cells.Item(17,2).Value2 = "Maxima"
$ws2.Cells.Item(17,3).Value2 = "2 lives, massive damage, mode change, 1st stage buff & debuffs &reflect dmg 2nd stage DOT & summons (carriers/cores every 40 seconds or smthing) (doesn't move all the way down but deal damage upon each death and if a stage is not killed fast enough it does game ending damage every 4mins or smthing."
$ws2.Cells.Item(18,1).Value2 = "[3-46]"
$ws2.Cells.Item(18,2).Value2 = "Ernesto"
$ws2.Cells.Item(18,3).Value2 = "slowly reduces max hp? Possible? If not hit for a time it accrues damage and hits with a burst on next hit. Enemies that detonate on fortress don't deal damage but ernesto holds their damage until it is big and then releases at once."
$ws2.Cells.Item(19,1).Value2 = "???"
$ws2.Cells.Item(19,2).Value2 = "??? (secret stage boss)"
$ws2.Cells.Item(19,3).Value2 = "3 lives, invulnerable during revives. Every stage it picks random abilities among the following, (teleporting movement, spawner, DOT, Buff, Debuff, Taunt & pull, bullets turn invisible, enemies turn invisible periodically) with the numbers chosen as 2,4,6`nduring reviving it stops production of bullets & cooldowns`nduring its second reviving it summons bunch of other stuff."

# 4) Sheet2 view: wider selection now that the table has grown.
$ws2.Activate()
$ws2.Range("I32").Select()

# 5) Sheet1 ("Basic enemies") view: scroll back to the top and move the selection.
$ws1.Activate()
$ws1.Range("L35").Select()

